$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The replacement figures include plain numeric-looking strings (e.g. "0.17")
# that must stay as literal TEXT (shared-string) cells, matching the rest of
# the table, instead of being auto-coerced to numbers by Excel's input
# parsing. Temporarily mark the range as Text so every value below is stored
# verbatim, then restore the original (General/"Normal") formatting so the
# cells end up styled exactly as they were before the edit.
$dataRange = $ws.Range("B2:D4")
$dataRange.NumberFormat = "@"

$ws.Range("B2").Value = "0.17"
$ws.Range("B3").Value = "-0.01"
$ws.Range("B4").Value = "-0.09"

$ws.Range("C2").Value = "44.29***"
$ws.Range("C3").Value = "2.21***"
$ws.Range("C4").Value = "0.98"

$ws.Range("D2").Value = "-0.89"
$ws.Range("D3").Value = "0.46***"
$ws.Range("D4").Value = "0.82*"

$dataRange.Style = "Normal"
